# Bug investments + updated lead times
# Update the "emission_cap" sheet: revise the 2040 and 2050 cap percentages
# (lead-time bug fix) and move the active selection to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("emission_cap")

# Row 5 (Year 2040): Percentage 30 -> 40
$ws.Range("B5").Value = 40

# Row 6 (Year 2050): Percentage 10 -> 20
$ws.Range("B6").Value = 20

# Leave the selection on B7, as in the saved workbook state
$ws.Range("B7").Select()
